$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)  # TeacherFreeSlot

# Update the time-slot strings (on the TeacherFreeSlot sheet) to the new
# dotted-time format, and fold the "8:30am-1:00pm;2:00pm-5:00pm" class
# slots into the same "8.30-11.30;14.00-17.00" slot used elsewhere.
$ws3.Range("C2").Value = "8.30-11.30;14.00-17.00"
$ws3.Range("E2").Value = "8.30-11.30;14.00-17.00"
$ws3.Range("F2").Value = "9.00-12.00;14.00-17.00"
$ws3.Range("G2").Value = "9.00-12.00;14.00-17.00"

$ws3.Range("C3").Value = "8.30-11.30;14.00-17.00"
$ws3.Range("E3").Value = "8.30-11.30;14.00-17.00"
$ws3.Range("G3").Value = "14.00-17.00"
$ws3.Range("D3").Value = "8.30-13.00"
$ws3.Range("F3").Value = "8.30-13.00"

# Make TeacherFreeSlot the active sheet/tab with D3 selected, matching the
# workbook's new activeTab and the sheet's new selection.
$ws3.Activate()
$ws3.Range("D3").Select()
